# Fill in the "Alias", "Estado" and "Esfuerzo (hrs)" columns for the first
# few use cases on the "Casos de Uso" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# Row 5 - CU-01: type the Alias / Estado / Esfuerzo directly.
$ws.Range("D5").Value = "JHAIR GARCIA CEBALLOS"
$ws.Range("E5").Value = "EN PROCESO"
$ws.Range("F5").Value = 8

# Copy the Estado/Esfuerzo (E5:F5) down onto rows 6-9 one row at a time.
$ws.Range("E5:F5").Copy()
$ws.Range("E6:F6").PasteSpecial()
$ws.Range("E5:F5").Copy()
$ws.Range("E7:F7").PasteSpecial()
$ws.Range("E5:F5").Copy()
$ws.Range("E8:F8").PasteSpecial()
$ws.Range("E5:F5").Copy()
$ws.Range("E9:F9").PasteSpecial()

# Row 10 - CU-06 (only the Alias is filled in, typed before the remaining
# Alias values below)
$ws.Range("D10").Value = "YAREL"

# Row 6 - CU-02
$ws.Range("D6").Value = "RUBEN JAHIR  RIVERA M."

# Row 8 - CU-04
$ws.Range("D8").Value = "ANDRES FLORES SANTAMARIA"

# Row 7 - CU-03
$ws.Range("D7").Value = "ANDRE HERNANDEZ C."

# Row 9 - CU-05
$ws.Range("D9").Value = "YAREL BAIZABAL VARGAS"

$ws.Range("D10").Select()
